$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.497.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.197.64'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.62%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.51'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.86%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.194.89'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.52%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.13%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.73%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.88%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.12%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.722.95'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.470.28'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.201.21'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.77'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.73%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.11'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.08'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.59%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.37%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +16.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.94'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.30'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.15%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.07%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.54'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '512.79'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.87'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0901'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0425'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.54%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.124'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.95%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.72%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.84%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +17.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.925.40'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.62%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.11%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.02%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +10.91%  '
